$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the 9 new match rows (263-271) for matchdays 19/04/13, 20/04/13, 21/04/13
# Row 263
$ws.Cells.Item(263, 1).Value = "D1"
$ws.Cells.Item(263, 2).Value = "19/04/13"
$ws.Cells.Item(263, 3).Value = "M'gladbach"
$ws.Cells.Item(263, 4).Value = "Augsburg"
$ws.Cells.Item(263, 5).Value = 1
$ws.Cells.Item(263, 6).Value = 0
$ws.Cells.Item(263, 7).Value = "H"
$ws.Cells.Item(263, 8).Value = 1
$ws.Cells.Item(263, 9).Value = 0
$ws.Cells.Item(263, 10).Value = "H"
$ws.Cells.Item(263, 11).Value = 23
$ws.Cells.Item(263, 12).Value = 11
$ws.Cells.Item(263, 13).Value = 5
$ws.Cells.Item(263, 14).Value = 2
$ws.Cells.Item(263, 15).Value = 18
$ws.Cells.Item(263, 16).Value = 14
$ws.Cells.Item(263, 17).Value = 6
$ws.Cells.Item(263, 18).Value = 5
$ws.Cells.Item(263, 19).Value = 4
$ws.Cells.Item(263, 20).Value = 1
$ws.Cells.Item(263, 21).Value = 0
$ws.Cells.Item(263, 22).Value = 1
$ws.Cells.Item(263, 23).Value = 2.3
$ws.Cells.Item(263, 24).Value = 3.2
$ws.Cells.Item(263, 25).Value = 3.1
$ws.Cells.Item(263, 26).Value = 2.35
$ws.Cells.Item(263, 27).Value = 3.25
$ws.Cells.Item(263, 28).Value = 3
$ws.Cells.Item(263, 29).Value = 2.35
$ws.Cells.Item(263, 30).Value = 3.25
$ws.Cells.Item(263, 31).Value = 3
$ws.Cells.Item(263, 32).Value = 2.1
$ws.Cells.Item(263, 33).Value = 3.3
$ws.Cells.Item(263, 34).Value = 3.3
$ws.Cells.Item(263, 35).Value = 2.25
$ws.Cells.Item(263, 36).Value = 3.3
$ws.Cells.Item(263, 37).Value = 3.2
$ws.Cells.Item(263, 38).Value = 2.39
$ws.Cells.Item(263, 39).Value = 3.41
$ws.Cells.Item(263, 40).Value = 3.24
$ws.Cells.Item(263, 41).Value = 2.3
$ws.Cells.Item(263, 42).Value = 3.3
$ws.Cells.Item(263, 43).Value = 3.1
$ws.Cells.Item(263, 44).Value = 2.3
$ws.Cells.Item(263, 45).Value = 3.2
$ws.Cells.Item(263, 46).Value = 3
$ws.Cells.Item(263, 47).Value = 2.45
$ws.Cells.Item(263, 48).Value = 3.4
$ws.Cells.Item(263, 49).Value = 3.12
$ws.Cells.Item(263, 50).Value = 2.4
$ws.Cells.Item(263, 51).Value = 3.3
$ws.Cells.Item(263, 52).Value = 2.9
$ws.Cells.Item(263, 53).Value = 38
$ws.Cells.Item(263, 54).Value = 2.41
$ws.Cells.Item(263, 55).Value = 2.31
$ws.Cells.Item(263, 56).Value = 3.41
$ws.Cells.Item(263, 57).Value = 3.28
$ws.Cells.Item(263, 58).Value = 3.35
$ws.Cells.Item(263, 59).Value = 3.15
$ws.Cells.Item(263, 60).Value = 35
$ws.Cells.Item(263, 61).Value = 2.13
$ws.Cells.Item(263, 62).Value = 2.05
$ws.Cells.Item(263, 63).Value = 1.83
$ws.Cells.Item(263, 64).Value = 1.76
$ws.Cells.Item(263, 65).Value = 22
$ws.Cells.Item(263, 66).Value = 0
$ws.Cells.Item(263, 67).Value = 1.72
$ws.Cells.Item(263, 68).Value = 1.67
$ws.Cells.Item(263, 69).Value = 2.32
$ws.Cells.Item(263, 70).Value = 2.26

# Row 264
$ws.Cells.Item(264, 1).Value = "D1"
$ws.Cells.Item(264, 2).Value = "20/04/13"
$ws.Cells.Item(264, 3).Value = "Dortmund"
$ws.Cells.Item(264, 4).Value = "Mainz"
$ws.Cells.Item(264, 5).Value = 2
$ws.Cells.Item(264, 6).Value = 0
$ws.Cells.Item(264, 7).Value = "H"
$ws.Cells.Item(264, 8).Value = 1
$ws.Cells.Item(264, 9).Value = 0
$ws.Cells.Item(264, 10).Value = "H"
$ws.Cells.Item(264, 11).Value = 22
$ws.Cells.Item(264, 12).Value = 5
$ws.Cells.Item(264, 13).Value = 7
$ws.Cells.Item(264, 14).Value = 0
$ws.Cells.Item(264, 15).Value = 12
$ws.Cells.Item(264, 16).Value = 16
$ws.Cells.Item(264, 17).Value = 8
$ws.Cells.Item(264, 18).Value = 4
$ws.Cells.Item(264, 19).Value = 3
$ws.Cells.Item(264, 20).Value = 3
$ws.Cells.Item(264, 21).Value = 0
$ws.Cells.Item(264, 22).Value = 0
$ws.Cells.Item(264, 23).Value = 1.44
$ws.Cells.Item(264, 24).Value = 4.5
$ws.Cells.Item(264, 25).Value = 6.5
$ws.Cells.Item(264, 26).Value = 1.45
$ws.Cells.Item(264, 27).Value = 4.5
$ws.Cells.Item(264, 28).Value = 6.5
$ws.Cells.Item(264, 29).Value = 1.45
$ws.Cells.Item(264, 30).Value = 4.5
$ws.Cells.Item(264, 31).Value = 6.5
$ws.Cells.Item(264, 32).Value = 1.5
$ws.Cells.Item(264, 33).Value = 4
$ws.Cells.Item(264, 34).Value = 6.1
$ws.Cells.Item(264, 35).Value = 1.5
$ws.Cells.Item(264, 36).Value = 4.33
$ws.Cells.Item(264, 37).Value = 6
$ws.Cells.Item(264, 38).Value = 1.52
$ws.Cells.Item(264, 39).Value = 4.68
$ws.Cells.Item(264, 40).Value = 6.77
$ws.Cells.Item(264, 41).Value = 1.53
$ws.Cells.Item(264, 42).Value = 4
$ws.Cells.Item(264, 43).Value = 6
$ws.Cells.Item(264, 44).Value = 1.44
$ws.Cells.Item(264, 45).Value = 4.33
$ws.Cells.Item(264, 46).Value = 6.25
$ws.Cells.Item(264, 47).Value = 1.53
$ws.Cells.Item(264, 48).Value = 4.6
$ws.Cells.Item(264, 49).Value = 6.5
$ws.Cells.Item(264, 50).Value = 1.44
$ws.Cells.Item(264, 51).Value = 4.33
$ws.Cells.Item(264, 52).Value = 7
$ws.Cells.Item(264, 53).Value = 38
$ws.Cells.Item(264, 54).Value = 1.53
$ws.Cells.Item(264, 55).Value = 1.48
$ws.Cells.Item(264, 56).Value = 4.68
$ws.Cells.Item(264, 57).Value = 4.3
$ws.Cells.Item(264, 58).Value = 7.1
$ws.Cells.Item(264, 59).Value = 6.45
$ws.Cells.Item(264, 60).Value = 26
$ws.Cells.Item(264, 61).Value = 1.65
$ws.Cells.Item(264, 62).Value = 1.57
$ws.Cells.Item(264, 63).Value = 2.5
$ws.Cells.Item(264, 64).Value = 2.31
$ws.Cells.Item(264, 65).Value = 24
$ws.Cells.Item(264, 66).Value = -1
$ws.Cells.Item(264, 67).Value = 1.85
$ws.Cells.Item(264, 68).Value = 1.79
$ws.Cells.Item(264, 69).Value = 2.15
$ws.Cells.Item(264, 70).Value = 2.08

# Row 265
$ws.Cells.Item(265, 1).Value = "D1"
$ws.Cells.Item(265, 2).Value = "20/04/13"
$ws.Cells.Item(265, 3).Value = "Ein Frankfurt"
$ws.Cells.Item(265, 4).Value = "Schalke 04"
$ws.Cells.Item(265, 5).Value = 1
$ws.Cells.Item(265, 6).Value = 0
$ws.Cells.Item(265, 7).Value = "H"
$ws.Cells.Item(265, 8).Value = 1
$ws.Cells.Item(265, 9).Value = 0
$ws.Cells.Item(265, 10).Value = "H"
$ws.Cells.Item(265, 11).Value = 10
$ws.Cells.Item(265, 12).Value = 10
$ws.Cells.Item(265, 13).Value = 3
$ws.Cells.Item(265, 14).Value = 7
$ws.Cells.Item(265, 15).Value = 16
$ws.Cells.Item(265, 16).Value = 27
$ws.Cells.Item(265, 17).Value = 0
$ws.Cells.Item(265, 18).Value = 4
$ws.Cells.Item(265, 19).Value = 3
$ws.Cells.Item(265, 20).Value = 3
$ws.Cells.Item(265, 21).Value = 0
$ws.Cells.Item(265, 22).Value = 0
$ws.Cells.Item(265, 23).Value = 3
$ws.Cells.Item(265, 24).Value = 3.4
$ws.Cells.Item(265, 25).Value = 2.3
$ws.Cells.Item(265, 26).Value = 2.95
$ws.Cells.Item(265, 27).Value = 3.4
$ws.Cells.Item(265, 28).Value = 2.3
$ws.Cells.Item(265, 29).Value = 2.95
$ws.Cells.Item(265, 30).Value = 3.4
$ws.Cells.Item(265, 31).Value = 2.3
$ws.Cells.Item(265, 32).Value = 2.9
$ws.Cells.Item(265, 33).Value = 3.3
$ws.Cells.Item(265, 34).Value = 2.3
$ws.Cells.Item(265, 35).Value = 3
$ws.Cells.Item(265, 36).Value = 3.4
$ws.Cells.Item(265, 37).Value = 2.29
$ws.Cells.Item(265, 38).Value = 3.19
$ws.Cells.Item(265, 39).Value = 3.55
$ws.Cells.Item(265, 40).Value = 2.35
$ws.Cells.Item(265, 41).Value = 3.1
$ws.Cells.Item(265, 42).Value = 3.3
$ws.Cells.Item(265, 43).Value = 2.3
$ws.Cells.Item(265, 44).Value = 3
$ws.Cells.Item(265, 45).Value = 3.3
$ws.Cells.Item(265, 46).Value = 2.25
$ws.Cells.Item(265, 47).Value = 3
$ws.Cells.Item(265, 48).Value = 3.5
$ws.Cells.Item(265, 49).Value = 2.4
$ws.Cells.Item(265, 50).Value = 3.1
$ws.Cells.Item(265, 51).Value = 3.4
$ws.Cells.Item(265, 52).Value = 2.25
$ws.Cells.Item(265, 53).Value = 38
$ws.Cells.Item(265, 54).Value = 3.2
$ws.Cells.Item(265, 55).Value = 3.01
$ws.Cells.Item(265, 56).Value = 3.6
$ws.Cells.Item(265, 57).Value = 3.4
$ws.Cells.Item(265, 58).Value = 2.4
$ws.Cells.Item(265, 59).Value = 2.3
$ws.Cells.Item(265, 60).Value = 35
$ws.Cells.Item(265, 61).Value = 1.8
$ws.Cells.Item(265, 62).Value = 1.72
$ws.Cells.Item(265, 63).Value = 2.22
$ws.Cells.Item(265, 64).Value = 2.09
$ws.Cells.Item(265, 65).Value = 23
$ws.Cells.Item(265, 66).Value = 0.25
$ws.Cells.Item(265, 67).Value = 1.91
$ws.Cells.Item(265, 68).Value = 1.88
$ws.Cells.Item(265, 69).Value = 2.04
$ws.Cells.Item(265, 70).Value = 2

# Row 266
$ws.Cells.Item(266, 1).Value = "D1"
$ws.Cells.Item(266, 2).Value = "20/04/13"
$ws.Cells.Item(266, 3).Value = "Hamburg"
$ws.Cells.Item(266, 4).Value = "Fortuna Dusseldorf"
$ws.Cells.Item(266, 5).Value = 2
$ws.Cells.Item(266, 6).Value = 1
$ws.Cells.Item(266, 7).Value = "H"
$ws.Cells.Item(266, 8).Value = 2
$ws.Cells.Item(266, 9).Value = 1
$ws.Cells.Item(266, 10).Value = "H"
$ws.Cells.Item(266, 11).Value = 17
$ws.Cells.Item(266, 12).Value = 12
$ws.Cells.Item(266, 13).Value = 4
$ws.Cells.Item(266, 14).Value = 5
$ws.Cells.Item(266, 15).Value = 21
$ws.Cells.Item(266, 16).Value = 19
$ws.Cells.Item(266, 17).Value = 5
$ws.Cells.Item(266, 18).Value = 3
$ws.Cells.Item(266, 19).Value = 2
$ws.Cells.Item(266, 20).Value = 3
$ws.Cells.Item(266, 21).Value = 0
$ws.Cells.Item(266, 22).Value = 0
$ws.Cells.Item(266, 23).Value = 1.83
$ws.Cells.Item(266, 24).Value = 3.5
$ws.Cells.Item(266, 25).Value = 4.33
$ws.Cells.Item(266, 26).Value = 1.83
$ws.Cells.Item(266, 27).Value = 3.6
$ws.Cells.Item(266, 28).Value = 4.2
$ws.Cells.Item(266, 29).Value = 1.83
$ws.Cells.Item(266, 30).Value = 3.6
$ws.Cells.Item(266, 31).Value = 4.2
$ws.Cells.Item(266, 32).Value = 1.75
$ws.Cells.Item(266, 33).Value = 3.7
$ws.Cells.Item(266, 34).Value = 4.2
$ws.Cells.Item(266, 35).Value = 1.83
$ws.Cells.Item(266, 36).Value = 3.5
$ws.Cells.Item(266, 37).Value = 4.33
$ws.Cells.Item(266, 38).Value = 1.83
$ws.Cells.Item(266, 39).Value = 3.78
$ws.Cells.Item(266, 40).Value = 4.79
$ws.Cells.Item(266, 41).Value = 1.85
$ws.Cells.Item(266, 42).Value = 3.3
$ws.Cells.Item(266, 43).Value = 4.5
$ws.Cells.Item(266, 44).Value = 1.8
$ws.Cells.Item(266, 45).Value = 3.4
$ws.Cells.Item(266, 46).Value = 4.33
$ws.Cells.Item(266, 47).Value = 1.85
$ws.Cells.Item(266, 48).Value = 3.75
$ws.Cells.Item(266, 49).Value = 4.6
$ws.Cells.Item(266, 50).Value = 1.83
$ws.Cells.Item(266, 51).Value = 3.5
$ws.Cells.Item(266, 52).Value = 4.33
$ws.Cells.Item(266, 53).Value = 38
$ws.Cells.Item(266, 54).Value = 1.9
$ws.Cells.Item(266, 55).Value = 1.82
$ws.Cells.Item(266, 56).Value = 3.76
$ws.Cells.Item(266, 57).Value = 3.58
$ws.Cells.Item(266, 58).Value = 5
$ws.Cells.Item(266, 59).Value = 4.32
$ws.Cells.Item(266, 60).Value = 35
$ws.Cells.Item(266, 61).Value = 2.02
$ws.Cells.Item(266, 62).Value = 1.86
$ws.Cells.Item(266, 63).Value = 2.09
$ws.Cells.Item(266, 64).Value = 1.93
$ws.Cells.Item(266, 65).Value = 23
$ws.Cells.Item(266, 66).Value = -0.75
$ws.Cells.Item(266, 67).Value = 2.13
$ws.Cells.Item(266, 68).Value = 2.09
$ws.Cells.Item(266, 69).Value = 1.86
$ws.Cells.Item(266, 70).Value = 1.79

# Row 267
$ws.Cells.Item(267, 1).Value = "D1"
$ws.Cells.Item(267, 2).Value = "20/04/13"
$ws.Cells.Item(267, 3).Value = "Hannover"
$ws.Cells.Item(267, 4).Value = "Bayern Munich"
$ws.Cells.Item(267, 5).Value = 1
$ws.Cells.Item(267, 6).Value = 6
$ws.Cells.Item(267, 7).Value = "A"
$ws.Cells.Item(267, 8).Value = 0
$ws.Cells.Item(267, 9).Value = 3
$ws.Cells.Item(267, 10).Value = "A"
$ws.Cells.Item(267, 11).Value = 7
$ws.Cells.Item(267, 12).Value = 14
$ws.Cells.Item(267, 13).Value = 4
$ws.Cells.Item(267, 14).Value = 5
$ws.Cells.Item(267, 15).Value = 7
$ws.Cells.Item(267, 16).Value = 9
$ws.Cells.Item(267, 17).Value = 2
$ws.Cells.Item(267, 18).Value = 9
$ws.Cells.Item(267, 19).Value = 0
$ws.Cells.Item(267, 20).Value = 1
$ws.Cells.Item(267, 21).Value = 0
$ws.Cells.Item(267, 22).Value = 0
$ws.Cells.Item(267, 23).Value = 5.5
$ws.Cells.Item(267, 24).Value = 4
$ws.Cells.Item(267, 25).Value = 1.57
$ws.Cells.Item(267, 26).Value = 5.75
$ws.Cells.Item(267, 27).Value = 3.9
$ws.Cells.Item(267, 28).Value = 1.57
$ws.Cells.Item(267, 29).Value = 5.75
$ws.Cells.Item(267, 30).Value = 3.9
$ws.Cells.Item(267, 31).Value = 1.57
$ws.Cells.Item(267, 32).Value = 4.3
$ws.Cells.Item(267, 33).Value = 3.9
$ws.Cells.Item(267, 34).Value = 1.7
$ws.Cells.Item(267, 35).Value = 5.5
$ws.Cells.Item(267, 36).Value = 4
$ws.Cells.Item(267, 37).Value = 1.57
$ws.Cells.Item(267, 38).Value = 6.37
$ws.Cells.Item(267, 39).Value = 4.25
$ws.Cells.Item(267, 40).Value = 1.59
$ws.Cells.Item(267, 41).Value = 5.5
$ws.Cells.Item(267, 42).Value = 3.75
$ws.Cells.Item(267, 43).Value = 1.62
$ws.Cells.Item(267, 44).Value = 6
$ws.Cells.Item(267, 45).Value = 3.75
$ws.Cells.Item(267, 46).Value = 1.53
$ws.Cells.Item(267, 47).Value = 5.75
$ws.Cells.Item(267, 48).Value = 4.3
$ws.Cells.Item(267, 49).Value = 1.62
$ws.Cells.Item(267, 50).Value = 5.5
$ws.Cells.Item(267, 51).Value = 3.8
$ws.Cells.Item(267, 52).Value = 1.62
$ws.Cells.Item(267, 53).Value = 38
$ws.Cells.Item(267, 54).Value = 6.4
$ws.Cells.Item(267, 55).Value = 5.67
$ws.Cells.Item(267, 56).Value = 4.31
$ws.Cells.Item(267, 57).Value = 3.93
$ws.Cells.Item(267, 58).Value = 1.7
$ws.Cells.Item(267, 59).Value = 1.59
$ws.Cells.Item(267, 60).Value = 35
$ws.Cells.Item(267, 61).Value = 1.77
$ws.Cells.Item(267, 62).Value = 1.7
$ws.Cells.Item(267, 63).Value = 2.26
$ws.Cells.Item(267, 64).Value = 2.12
$ws.Cells.Item(267, 65).Value = 24
$ws.Cells.Item(267, 66).Value = 1
$ws.Cells.Item(267, 67).Value = 1.95
$ws.Cells.Item(267, 68).Value = 1.88
$ws.Cells.Item(267, 69).Value = 2.03
$ws.Cells.Item(267, 70).Value = 1.98

# Row 268
$ws.Cells.Item(268, 1).Value = "D1"
$ws.Cells.Item(268, 2).Value = "20/04/13"
$ws.Cells.Item(268, 3).Value = "Leverkusen"
$ws.Cells.Item(268, 4).Value = "Hoffenheim"
$ws.Cells.Item(268, 5).Value = 5
$ws.Cells.Item(268, 6).Value = 0
$ws.Cells.Item(268, 7).Value = "H"
$ws.Cells.Item(268, 8).Value = 2
$ws.Cells.Item(268, 9).Value = 0
$ws.Cells.Item(268, 10).Value = "H"
$ws.Cells.Item(268, 11).Value = 22
$ws.Cells.Item(268, 12).Value = 5
$ws.Cells.Item(268, 13).Value = 9
$ws.Cells.Item(268, 14).Value = 2
$ws.Cells.Item(268, 15).Value = 11
$ws.Cells.Item(268, 16).Value = 9
$ws.Cells.Item(268, 17).Value = 5
$ws.Cells.Item(268, 18).Value = 1
$ws.Cells.Item(268, 19).Value = 1
$ws.Cells.Item(268, 20).Value = 0
$ws.Cells.Item(268, 21).Value = 0
$ws.Cells.Item(268, 22).Value = 1
$ws.Cells.Item(268, 23).Value = 1.44
$ws.Cells.Item(268, 24).Value = 4.5
$ws.Cells.Item(268, 25).Value = 6.5
$ws.Cells.Item(268, 26).Value = 1.45
$ws.Cells.Item(268, 27).Value = 4.75
$ws.Cells.Item(268, 28).Value = 6
$ws.Cells.Item(268, 29).Value = 1.45
$ws.Cells.Item(268, 30).Value = 4.75
$ws.Cells.Item(268, 31).Value = 6
$ws.Cells.Item(268, 32).Value = 1.45
$ws.Cells.Item(268, 33).Value = 4
$ws.Cells.Item(268, 34).Value = 7
$ws.Cells.Item(268, 35).Value = 1.44
$ws.Cells.Item(268, 36).Value = 4.5
$ws.Cells.Item(268, 37).Value = 6.5
$ws.Cells.Item(268, 38).Value = 1.49
$ws.Cells.Item(268, 39).Value = 4.76
$ws.Cells.Item(268, 40).Value = 7.31
$ws.Cells.Item(268, 41).Value = 1.5
$ws.Cells.Item(268, 42).Value = 4
$ws.Cells.Item(268, 43).Value = 7
$ws.Cells.Item(268, 44).Value = 1.44
$ws.Cells.Item(268, 45).Value = 4
$ws.Cells.Item(268, 46).Value = 7
$ws.Cells.Item(268, 47).Value = 1.5
$ws.Cells.Item(268, 48).Value = 4.6
$ws.Cells.Item(268, 49).Value = 7
$ws.Cells.Item(268, 50).Value = 1.44
$ws.Cells.Item(268, 51).Value = 4.33
$ws.Cells.Item(268, 52).Value = 7
$ws.Cells.Item(268, 53).Value = 38
$ws.Cells.Item(268, 54).Value = 1.5
$ws.Cells.Item(268, 55).Value = 1.47
$ws.Cells.Item(268, 56).Value = 4.8
$ws.Cells.Item(268, 57).Value = 4.31
$ws.Cells.Item(268, 58).Value = 7.6
$ws.Cells.Item(268, 59).Value = 6.86
$ws.Cells.Item(268, 60).Value = 26
$ws.Cells.Item(268, 61).Value = 1.68
$ws.Cells.Item(268, 62).Value = 1.62
$ws.Cells.Item(268, 63).Value = 2.38
$ws.Cells.Item(268, 64).Value = 2.22
$ws.Cells.Item(268, 65).Value = 24
$ws.Cells.Item(268, 66).Value = -1
$ws.Cells.Item(268, 67).Value = 1.8
$ws.Cells.Item(268, 68).Value = 1.76
$ws.Cells.Item(268, 69).Value = 2.16
$ws.Cells.Item(268, 70).Value = 2.11

# Row 269
$ws.Cells.Item(269, 1).Value = "D1"
$ws.Cells.Item(269, 2).Value = "20/04/13"
$ws.Cells.Item(269, 3).Value = "Werder Bremen"
$ws.Cells.Item(269, 4).Value = "Wolfsburg"
$ws.Cells.Item(269, 5).Value = 0
$ws.Cells.Item(269, 6).Value = 3
$ws.Cells.Item(269, 7).Value = "A"
$ws.Cells.Item(269, 8).Value = 0
$ws.Cells.Item(269, 9).Value = 2
$ws.Cells.Item(269, 10).Value = "A"
$ws.Cells.Item(269, 11).Value = 8
$ws.Cells.Item(269, 12).Value = 8
$ws.Cells.Item(269, 13).Value = 3
$ws.Cells.Item(269, 14).Value = 5
$ws.Cells.Item(269, 15).Value = 18
$ws.Cells.Item(269, 16).Value = 16
$ws.Cells.Item(269, 17).Value = 6
$ws.Cells.Item(269, 18).Value = 3
$ws.Cells.Item(269, 19).Value = 3
$ws.Cells.Item(269, 20).Value = 1
$ws.Cells.Item(269, 21).Value = 0
$ws.Cells.Item(269, 22).Value = 0
$ws.Cells.Item(269, 23).Value = 2
$ws.Cells.Item(269, 24).Value = 3.6
$ws.Cells.Item(269, 25).Value = 3.5
$ws.Cells.Item(269, 26).Value = 2.05
$ws.Cells.Item(269, 27).Value = 3.4
$ws.Cells.Item(269, 28).Value = 3.5
$ws.Cells.Item(269, 29).Value = 2.05
$ws.Cells.Item(269, 30).Value = 3.4
$ws.Cells.Item(269, 31).Value = 3.5
$ws.Cells.Item(269, 32).Value = 2.1
$ws.Cells.Item(269, 33).Value = 3.3
$ws.Cells.Item(269, 34).Value = 3.3
$ws.Cells.Item(269, 35).Value = 2.1
$ws.Cells.Item(269, 36).Value = 3.4
$ws.Cells.Item(269, 37).Value = 3.4
$ws.Cells.Item(269, 38).Value = 2.1
$ws.Cells.Item(269, 39).Value = 3.69
$ws.Cells.Item(269, 40).Value = 3.66
$ws.Cells.Item(269, 41).Value = 2.15
$ws.Cells.Item(269, 42).Value = 3.5
$ws.Cells.Item(269, 43).Value = 3.2
$ws.Cells.Item(269, 44).Value = 2.05
$ws.Cells.Item(269, 45).Value = 3.4
$ws.Cells.Item(269, 46).Value = 3.3
$ws.Cells.Item(269, 47).Value = 2.15
$ws.Cells.Item(269, 48).Value = 3.5
$ws.Cells.Item(269, 49).Value = 3.6
$ws.Cells.Item(269, 50).Value = 2.15
$ws.Cells.Item(269, 51).Value = 3.4
$ws.Cells.Item(269, 52).Value = 3.3
$ws.Cells.Item(269, 53).Value = 38
$ws.Cells.Item(269, 54).Value = 2.15
$ws.Cells.Item(269, 55).Value = 2.09
$ws.Cells.Item(269, 56).Value = 3.75
$ws.Cells.Item(269, 57).Value = 3.48
$ws.Cells.Item(269, 58).Value = 3.59
$ws.Cells.Item(269, 59).Value = 3.38
$ws.Cells.Item(269, 60).Value = 27
$ws.Cells.Item(269, 61).Value = 1.68
$ws.Cells.Item(269, 62).Value = 1.59
$ws.Cells.Item(269, 63).Value = 2.48
$ws.Cells.Item(269, 64).Value = 2.26
$ws.Cells.Item(269, 65).Value = 21
$ws.Cells.Item(269, 66).Value = -0.25
$ws.Cells.Item(269, 67).Value = 1.85
$ws.Cells.Item(269, 68).Value = 1.8
$ws.Cells.Item(269, 69).Value = 2.14
$ws.Cells.Item(269, 70).Value = 2.08

# Row 270
$ws.Cells.Item(270, 1).Value = "D1"
$ws.Cells.Item(270, 2).Value = "21/04/13"
$ws.Cells.Item(270, 3).Value = "Nurnberg"
$ws.Cells.Item(270, 4).Value = "Greuther Furth"
$ws.Cells.Item(270, 5).Value = 0
$ws.Cells.Item(270, 6).Value = 1
$ws.Cells.Item(270, 7).Value = "A"
$ws.Cells.Item(270, 8).Value = 0
$ws.Cells.Item(270, 9).Value = 1
$ws.Cells.Item(270, 10).Value = "A"
$ws.Cells.Item(270, 11).Value = 23
$ws.Cells.Item(270, 12).Value = 15
$ws.Cells.Item(270, 13).Value = 8
$ws.Cells.Item(270, 14).Value = 3
$ws.Cells.Item(270, 15).Value = 12
$ws.Cells.Item(270, 16).Value = 25
$ws.Cells.Item(270, 17).Value = 14
$ws.Cells.Item(270, 18).Value = 4
$ws.Cells.Item(270, 19).Value = 4
$ws.Cells.Item(270, 20).Value = 5
$ws.Cells.Item(270, 21).Value = 0
$ws.Cells.Item(270, 22).Value = 0
$ws.Cells.Item(270, 23).Value = 1.75
$ws.Cells.Item(270, 24).Value = 3.5
$ws.Cells.Item(270, 25).Value = 4.75
$ws.Cells.Item(270, 26).Value = 1.7
$ws.Cells.Item(270, 27).Value = 3.6
$ws.Cells.Item(270, 28).Value = 5
$ws.Cells.Item(270, 29).Value = 1.7
$ws.Cells.Item(270, 30).Value = 3.6
$ws.Cells.Item(270, 31).Value = 5
$ws.Cells.Item(270, 32).Value = 1.7
$ws.Cells.Item(270, 33).Value = 3.7
$ws.Cells.Item(270, 34).Value = 4.5
$ws.Cells.Item(270, 35).Value = 1.72
$ws.Cells.Item(270, 36).Value = 3.5
$ws.Cells.Item(270, 37).Value = 5
$ws.Cells.Item(270, 38).Value = 1.77
$ws.Cells.Item(270, 39).Value = 3.73
$ws.Cells.Item(270, 40).Value = 5.37
$ws.Cells.Item(270, 41).Value = 1.75
$ws.Cells.Item(270, 42).Value = 3.4
$ws.Cells.Item(270, 43).Value = 5
$ws.Cells.Item(270, 44).Value = 1.67
$ws.Cells.Item(270, 45).Value = 3.6
$ws.Cells.Item(270, 46).Value = 4.8
$ws.Cells.Item(270, 47).Value = 1.75
$ws.Cells.Item(270, 48).Value = 3.7
$ws.Cells.Item(270, 49).Value = 5.4
$ws.Cells.Item(270, 50).Value = 1.67
$ws.Cells.Item(270, 51).Value = 3.75
$ws.Cells.Item(270, 52).Value = 5
$ws.Cells.Item(270, 53).Value = 38
$ws.Cells.Item(270, 54).Value = 1.8
$ws.Cells.Item(270, 55).Value = 1.72
$ws.Cells.Item(270, 56).Value = 3.85
$ws.Cells.Item(270, 57).Value = 3.59
$ws.Cells.Item(270, 58).Value = 5.5
$ws.Cells.Item(270, 59).Value = 4.95
$ws.Cells.Item(270, 60).Value = 36
$ws.Cells.Item(270, 61).Value = 2.14
$ws.Cells.Item(270, 62).Value = 2.04
$ws.Cells.Item(270, 63).Value = 1.83
$ws.Cells.Item(270, 64).Value = 1.76
$ws.Cells.Item(270, 65).Value = 23
$ws.Cells.Item(270, 66).Value = -1
$ws.Cells.Item(270, 67).Value = 2.43
$ws.Cells.Item(270, 68).Value = 2.34
$ws.Cells.Item(270, 69).Value = 1.68
$ws.Cells.Item(270, 70).Value = 1.63

# Row 271
$ws.Cells.Item(271, 1).Value = "D1"
$ws.Cells.Item(271, 2).Value = "21/04/13"
$ws.Cells.Item(271, 3).Value = "Stuttgart"
$ws.Cells.Item(271, 4).Value = "Freiburg"
$ws.Cells.Item(271, 5).Value = 2
$ws.Cells.Item(271, 6).Value = 1
$ws.Cells.Item(271, 7).Value = "H"
$ws.Cells.Item(271, 8).Value = 2
$ws.Cells.Item(271, 9).Value = 0
$ws.Cells.Item(271, 10).Value = "H"
$ws.Cells.Item(271, 11).Value = 15
$ws.Cells.Item(271, 12).Value = 11
$ws.Cells.Item(271, 13).Value = 7
$ws.Cells.Item(271, 14).Value = 3
$ws.Cells.Item(271, 15).Value = 14
$ws.Cells.Item(271, 16).Value = 15
$ws.Cells.Item(271, 17).Value = 3
$ws.Cells.Item(271, 18).Value = 6
$ws.Cells.Item(271, 19).Value = 0
$ws.Cells.Item(271, 20).Value = 3
$ws.Cells.Item(271, 21).Value = 0
$ws.Cells.Item(271, 22).Value = 0
$ws.Cells.Item(271, 23).Value = 2.3
$ws.Cells.Item(271, 24).Value = 3.4
$ws.Cells.Item(271, 25).Value = 3
$ws.Cells.Item(271, 26).Value = 2.4
$ws.Cells.Item(271, 27).Value = 3.4
$ws.Cells.Item(271, 28).Value = 2.85
$ws.Cells.Item(271, 29).Value = 2.4
$ws.Cells.Item(271, 30).Value = 3.4
$ws.Cells.Item(271, 31).Value = 2.85
$ws.Cells.Item(271, 32).Value = 2.4
$ws.Cells.Item(271, 33).Value = 3.3
$ws.Cells.Item(271, 34).Value = 2.75
$ws.Cells.Item(271, 35).Value = 2.37
$ws.Cells.Item(271, 36).Value = 3.4
$ws.Cells.Item(271, 37).Value = 2.87
$ws.Cells.Item(271, 38).Value = 2.45
$ws.Cells.Item(271, 39).Value = 3.51
$ws.Cells.Item(271, 40).Value = 3.06
$ws.Cells.Item(271, 41).Value = 2.5
$ws.Cells.Item(271, 42).Value = 3.2
$ws.Cells.Item(271, 43).Value = 2.8
$ws.Cells.Item(271, 44).Value = 2.3
$ws.Cells.Item(271, 45).Value = 3.3
$ws.Cells.Item(271, 46).Value = 2.88
$ws.Cells.Item(271, 47).Value = 2.5
$ws.Cells.Item(271, 48).Value = 3.5
$ws.Cells.Item(271, 49).Value = 3
$ws.Cells.Item(271, 50).Value = 2.38
$ws.Cells.Item(271, 51).Value = 3.4
$ws.Cells.Item(271, 52).Value = 2.88
$ws.Cells.Item(271, 53).Value = 38
$ws.Cells.Item(271, 54).Value = 2.55
$ws.Cells.Item(271, 55).Value = 2.4
$ws.Cells.Item(271, 56).Value = 3.51
$ws.Cells.Item(271, 57).Value = 3.33
$ws.Cells.Item(271, 58).Value = 3.13
$ws.Cells.Item(271, 59).Value = 2.92
$ws.Cells.Item(271, 60).Value = 34
$ws.Cells.Item(271, 61).Value = 1.94
$ws.Cells.Item(271, 62).Value = 1.87
$ws.Cells.Item(271, 63).Value = 2.02
$ws.Cells.Item(271, 64).Value = 1.92
$ws.Cells.Item(271, 65).Value = 21
$ws.Cells.Item(271, 66).Value = 0
$ws.Cells.Item(271, 67).Value = 1.76
$ws.Cells.Item(271, 68).Value = 1.72
$ws.Cells.Item(271, 69).Value = 2.26
$ws.Cells.Item(271, 70).Value = 2.17

# Update the selection to reflect the new active cell used while entering this data
$ws.Range("B269").Select()

